$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2425.6428
$ws.Range("I38").Value = 973.2222
$ws.Range("K38").Value = 2919.6666
$ws.Range("M38").Value = -2547.6666
$ws.Range("H40").Value = 3696.6667
$ws.Range("I40").Value = 4450
$ws.Range("J40").Value = 3546
$ws.Range("K40").Value = 4450
$ws.Range("L40").Value = 3546
$ws.Range("M40").Value = -4275
$ws.Range("N40").Value = -3896
$ws.Range("H132").Value = 2306.7551
$ws.Range("I132").Value = 2189.7954
$ws.Range("K132").Value = 6569.3862
$ws.Range("M132").Value = -4039.3862
$ws.Range("H137").Value = 21568
$ws.Range("I137").Value = 25250.334
$ws.Range("J137").Value = 4997.5
$ws.Range("K137").Value = 75751.00199999999
$ws.Range("L137").Value = 14992.5
$ws.Range("M137").Value = -73201.00199999999
$ws.Range("N137").Value = -20092.5
$ws.Range("H138").Value = 3959.7632
$ws.Range("I138").Value = 1782.375
$ws.Range("J138").Value = 4540.4
$ws.Range("K138").Value = 5347.125
$ws.Range("L138").Value = 13621.2
$ws.Range("M138").Value = -207.125
$ws.Range("N138").Value = -23901.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 49000
$ws.Range("J43").Value = 49000
$ws.Range("L43").Value = 49000
$ws.Range("N43").Value = -49626
$ws.Range("H122").Value = 1765880.4
$ws.Range("I122").Value = 4127.75
$ws.Range("J122").Value = 2594940.5
$ws.Range("K122").Value = 12383.25
$ws.Range("L122").Value = 7784821.5
$ws.Range("M122").Value = -9933.25
$ws.Range("N122").Value = -7789721.5
$ws.Range("H132").Value = 4049.4167
$ws.Range("I132").Value = 1986.1666
$ws.Range("K132").Value = 5958.4998
$ws.Range("M132").Value = -3428.4998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 6458.8
$ws.Range("I25").Value = 3766.3333
$ws.Range("J25").Value = 10497.5
$ws.Range("K25").Value = 3766.3333
$ws.Range("L25").Value = 10497.5
$ws.Range("M25").Value = -3531.3333
$ws.Range("N25").Value = -10967.5
$ws.Range("H86").Value = 8586.25
$ws.Range("I86").Value = 5929.1665
$ws.Range("J86").Value = 32500
$ws.Range("K86").Value = 5929.1665
$ws.Range("L86").Value = 32500
$ws.Range("M86").Value = -4806.1665
$ws.Range("N86").Value = -34746
$ws.Range("H89").Value = 8586.25
$ws.Range("I89").Value = 5929.1665
$ws.Range("J89").Value = 32500
$ws.Range("K89").Value = 29645.8325
$ws.Range("L89").Value = 162500
$ws.Range("M89").Value = -24029.8325
$ws.Range("N89").Value = -173732

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 5600
$ws.Range("I41").Value = 5600
$ws.Range("K41").Value = 5600
$ws.Range("M41").Value = -5172
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()
$ws.Range("H58").Value = 2314.0625
$ws.Range("I58").Value = 1681.4166
$ws.Range("K58").Value = 1681.4166
$ws.Range("M58").Value = -1478.4166
$ws.Range("H93").Value = 8200.5
$ws.Range("I93").Value = 4515
$ws.Range("K93").Value = 4515
$ws.Range("M93").Value = -2643
$ws.Range("H136").Value = 2314.0625
$ws.Range("I136").Value = 1681.4166
$ws.Range("K136").Value = 5044.2498
$ws.Range("M136").Value = -2494.2498

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 133.55556
$ws.Range("J2").Value = 209.8
$ws.Range("L2").Value = 1258.8
$ws.Range("N2").Value = -1484.8
$ws.Range("H6").Value = 78.875
$ws.Range("I6").Value = 78.875
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 236.625
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -123.625
$ws.Range("N6").ClearContents()
$ws.Range("H22").Value = 5456.6665
$ws.Range("I22").Value = 5435
$ws.Range("J22").Value = 5500
$ws.Range("K22").Value = 16305
$ws.Range("L22").Value = 16500
$ws.Range("M22").Value = -16136
$ws.Range("N22").Value = -16838
$ws.Range("H27").Value = 5456.6665
$ws.Range("I27").Value = 5435
$ws.Range("J27").Value = 5500
$ws.Range("K27").Value = 16305
$ws.Range("L27").Value = 16500
$ws.Range("M27").Value = -16203
$ws.Range("N27").Value = -16704
$ws.Range("H41").Value = 2822.2856
$ws.Range("I41").Value = 5244.5713
$ws.Range("J41").Value = 400
$ws.Range("K41").Value = 15733.7139
$ws.Range("L41").Value = 1200
$ws.Range("M41").Value = -15395.7139
$ws.Range("N41").Value = -1876
$ws.Range("H50").Value = 933.7273
$ws.Range("I50").Value = 518.375
$ws.Range("J50").Value = 2041.3334
$ws.Range("K50").Value = 1555.125
$ws.Range("L50").Value = 6124.0002
$ws.Range("M50").Value = -1074.125
$ws.Range("N50").Value = -7086.0002
$ws.Range("H53").Value = 933.7273
$ws.Range("I53").Value = 518.375
$ws.Range("J53").Value = 2041.3334
$ws.Range("K53").Value = 1555.125
$ws.Range("L53").Value = 6124.0002
$ws.Range("M53").Value = -1074.125
$ws.Range("N53").Value = -7086.0002
$ws.Range("H54").Value = 4699.75
$ws.Range("I54").Value = 933
$ws.Range("J54").Value = 16000
$ws.Range("K54").Value = 2799
$ws.Range("L54").Value = 48000
$ws.Range("M54").Value = -2240
$ws.Range("N54").Value = -49118
$ws.Range("H55").Value = 7704.3335
$ws.Range("I55").Value = 1124.6666
$ws.Range("K55").Value = 3373.9998
$ws.Range("M55").Value = -3196.9998
$ws.Range("H58").Value = 1401.6666
$ws.Range("I58").Value = 1401.6666
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 4204.9998
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -4076.9998
$ws.Range("N58").ClearContents()
$ws.Range("H61").Value = 244
$ws.Range("I61").Value = 212
$ws.Range("J61").Value = 324
$ws.Range("K61").Value = 636
$ws.Range("L61").Value = 972
$ws.Range("M61").Value = -421
$ws.Range("N61").Value = -1402
$ws.Range("H132").Value = 57684.555
$ws.Range("I132").Value = 888.5
$ws.Range("J132").Value = 73912
$ws.Range("K132").Value = 7996.5
$ws.Range("L132").Value = 665208
$ws.Range("M132").Value = -5466.5
$ws.Range("N132").Value = -670268
$ws.Range("H134").Value = 11071.454
$ws.Range("I134").Value = 11071.454
$ws.Range("K134").Value = 33214.362
$ws.Range("M134").Value = -28144.362
$ws.Range("H137").Value = 2492.9
$ws.Range("I137").Value = 2492.9
$ws.Range("K137").Value = 7478.700000000001
$ws.Range("M137").Value = -2378.700000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 34750
$ws.Range("J123").Value = 34750
$ws.Range("L123").Value = 34750
$ws.Range("N123").Value = -39650
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 20972.133
$ws.Range("J7").Value = 8745.916999999999
$ws.Range("L7").Value = 8745.916999999999
$ws.Range("N7").Value = -8969.916999999999
$ws.Range("H32").Value = 2877.5
$ws.Range("I32").Value = 2877.5
$ws.Range("K32").Value = 2877.5
$ws.Range("M32").Value = -2560.5
$ws.Range("H122").Value = 5914.5
$ws.Range("I122").Value = 3638
$ws.Range("J122").Value = 7540.5713
$ws.Range("K122").Value = 10914
$ws.Range("L122").Value = 22621.7139
$ws.Range("M122").Value = -8464
$ws.Range("N122").Value = -27521.7139
$ws.Range("H126").Value = 20972.133
$ws.Range("J126").Value = 8745.916999999999
$ws.Range("L126").Value = 26237.751
$ws.Range("N126").Value = -31177.751

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 13909
$ws.Range("I132").Value = 16097.704
$ws.Range("J132").Value = 7999.5
$ws.Range("K132").Value = 48293.112
$ws.Range("L132").Value = 23998.5
$ws.Range("M132").Value = -45763.112
$ws.Range("N132").Value = -29058.5
